$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.514.75'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").Value = '1.625.49'
$ws.Range("E3").Value = '  +0.04%  '
$c = $ws.Range("D4")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = $origStyle
$ws.Range("E4").Value = '  +0.01%  '
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '213.80'
$c.Style = $origStyle
$ws.Range("E5").Value = '  -0.38%  '
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.502'
$c.Style = $origStyle
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("E8").Value = '  -0.31%  '
$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0610'
$c.Style = $origStyle
$c = $ws.Range("D10")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '19.22'
$c.Style = $origStyle
$ws.Range("E10").Value = '  -0.93%  '
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("D12").Value = '1.852.36'
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").Value = '1.617.18'
$ws.Range("E13").Value = '  -0.60%  '
$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.05'
$c.Style = $origStyle
$ws.Range("E14").Value = '  -0.37%  '
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.512'
$c.Style = $origStyle
$ws.Range("E15").Value = '  -0.31%  '
$c = $ws.Range("D16")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '63.98'
$c.Style = $origStyle
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D17")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '234.72'
$c.Style = $origStyle
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '26.518.97'
$ws.Range("E18").Value = '  -0.86%  '
$c = $ws.Range("D19")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.75'
$c.Style = $origStyle
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").Value = '0.0₃0726'
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("E21").Value = '  +0.00%  '
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.32'
$c.Style = $origStyle
$ws.Range("E22").Value = '  -2.07%  '
$ws.Range("E23").Value = '  -1.81%  '
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.15'
$c.Style = $origStyle
$ws.Range("E24").Value = '  +0.24%  '
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '145.96'
$c.Style = $origStyle
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("E26").Value = '  +0.00%  '
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.08'
$c.Style = $origStyle
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  -0.68%  '
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '15.64'
$c.Style = $origStyle
$ws.Range("E29").Value = '  -0.35%  '
$c = $ws.Range("D30")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0493'
$c.Style = $origStyle
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("E31").Value = '  -0.35%  '
$ws.Range("D32").Value = '1.522.06'
$ws.Range("E32").Value = '  +3.27%  '
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.27'
$c.Style = $origStyle
$ws.Range("E33").Value = '  +0.24%  '
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.00'
$c.Style = $origStyle
$ws.Range("E34").Value = '  +0.01%  '
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.53'
$c.Style = $origStyle
$ws.Range("E35").Value = '  +2.38%  '
$ws.Range("E36").Value = '  -0.39%  '
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.572'
$c.Style = $origStyle
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("E38").Value = '  -0.81%  '
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.836'
$c.Style = $origStyle
$ws.Range("E39").Value = '  -0.52%  '
$ws.Range("E40").Value = '  -1.69%  '
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.21'
$c.Style = $origStyle
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").Value = '1.764.31'
$ws.Range("E43").Value = '  -0.01%  '
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '62.64'
$c.Style = $origStyle
$ws.Range("E44").Value = '  +0.73%  '
$c = $ws.Range("D45")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.761'
$c.Style = $origStyle
$ws.Range("E45").Value = '  -0.76%  '
$c = $ws.Range("D46")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.909'
$c.Style = $origStyle
$ws.Range("E46").Value = '  -5.01%  '
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '89.75'
$c.Style = $origStyle
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0102'
$ws.Range("E49").Value = '  +11.30%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D50")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0502'
$c.Style = $origStyle
$ws.Range("E50").Value = '  -0.72%  '
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0963'
$c.Style = $origStyle
$ws.Range("E51").Value = '  -0.30%  '
